$d = $word.ActiveDocument

# Update the date paragraph (first paragraph, outside the table)
$d.Content.Find.Execute("2023-07-21 Friday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-07-22 Saturday", 2) | Out-Null

# Update each answer cell in the 20x5 table, addressed by (row, col) to
# avoid ambiguity from duplicate old values appearing more than once.
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "59+0=59"
$t.Cell(1, 2).Range.Text = "9+55=64"
$t.Cell(1, 3).Range.Text = "53+13=66"
$t.Cell(1, 4).Range.Text = "54-9=45"
$t.Cell(1, 5).Range.Text = "45-16=29"

$t.Cell(2, 1).Range.Text = "49+32=81"
$t.Cell(2, 2).Range.Text = "17-10=7"
$t.Cell(2, 3).Range.Text = "96-66=30"
$t.Cell(2, 4).Range.Text = "53-40=13"
$t.Cell(2, 5).Range.Text = "2+4=6"

$t.Cell(3, 1).Range.Text = "40+36=76"
$t.Cell(3, 2).Range.Text = "8+45=53"
$t.Cell(3, 3).Range.Text = "72-48=24"
$t.Cell(3, 4).Range.Text = "94-85=9"
$t.Cell(3, 5).Range.Text = "90-1=89"

$t.Cell(4, 1).Range.Text = "41+39=80"
$t.Cell(4, 2).Range.Text = "34+33=67"
$t.Cell(4, 3).Range.Text = "2+72=74"
$t.Cell(4, 4).Range.Text = "7+24=31"
$t.Cell(4, 5).Range.Text = "79-15=64"

$t.Cell(5, 1).Range.Text = "67-2=65"
$t.Cell(5, 2).Range.Text = "95-22=73"
$t.Cell(5, 3).Range.Text = "40-18=22"
$t.Cell(5, 4).Range.Text = "0+91=91"
$t.Cell(5, 5).Range.Text = "73-37=36"

$t.Cell(6, 1).Range.Text = "68+15=83"
$t.Cell(6, 2).Range.Text = "36+8=44"
$t.Cell(6, 3).Range.Text = "26+72=98"
$t.Cell(6, 4).Range.Text = "21-8=13"
$t.Cell(6, 5).Range.Text = "52+26=78"

$t.Cell(7, 1).Range.Text = "69-17=52"
$t.Cell(7, 2).Range.Text = "23+60=83"
$t.Cell(7, 3).Range.Text = "29+39=68"
$t.Cell(7, 4).Range.Text = "18+58=76"
$t.Cell(7, 5).Range.Text = "98-23=75"

$t.Cell(8, 1).Range.Text = "20-16=4"
$t.Cell(8, 2).Range.Text = "78-45=33"
$t.Cell(8, 3).Range.Text = "99-54=45"
$t.Cell(8, 4).Range.Text = "46+22=68"
$t.Cell(8, 5).Range.Text = "75-37=38"

$t.Cell(9, 1).Range.Text = "20+58=78"
$t.Cell(9, 2).Range.Text = "73-24=49"
$t.Cell(9, 3).Range.Text = "47-39=8"
$t.Cell(9, 4).Range.Text = "89-73=16"
$t.Cell(9, 5).Range.Text = "84-65=19"

$t.Cell(10, 1).Range.Text = "2+12=14"
$t.Cell(10, 2).Range.Text = "73-44=29"
$t.Cell(10, 3).Range.Text = "91-15=76"
$t.Cell(10, 4).Range.Text = "56-43=13"
$t.Cell(10, 5).Range.Text = "87-51=36"

$t.Cell(11, 1).Range.Text = "8-8=0"
$t.Cell(11, 2).Range.Text = "20+34=54"
$t.Cell(11, 3).Range.Text = "92-76=16"
$t.Cell(11, 4).Range.Text = "14-0=14"
$t.Cell(11, 5).Range.Text = "55+14=69"

$t.Cell(12, 1).Range.Text = "77-11=66"
$t.Cell(12, 2).Range.Text = "46+30=76"
$t.Cell(12, 3).Range.Text = "67-58=9"
$t.Cell(12, 4).Range.Text = "62+31=93"
$t.Cell(12, 5).Range.Text = "46-32=14"

$t.Cell(13, 1).Range.Text = "84-17=67"
$t.Cell(13, 2).Range.Text = "29-4=25"
$t.Cell(13, 3).Range.Text = "53+30=83"
$t.Cell(13, 4).Range.Text = "83-27=56"
$t.Cell(13, 5).Range.Text = "40-31=9"

$t.Cell(14, 1).Range.Text = "0+17=17"
$t.Cell(14, 2).Range.Text = "83-35=48"
$t.Cell(14, 3).Range.Text = "0+95=95"
$t.Cell(14, 4).Range.Text = "71-39=32"
$t.Cell(14, 5).Range.Text = "93-37=56"

$t.Cell(15, 1).Range.Text = "18+15=33"
$t.Cell(15, 2).Range.Text = "91-40=51"
$t.Cell(15, 3).Range.Text = "81+1=82"
$t.Cell(15, 4).Range.Text = "81-43=38"
$t.Cell(15, 5).Range.Text = "8+90=98"

$t.Cell(16, 1).Range.Text = "88+9=97"
$t.Cell(16, 2).Range.Text = "31+63=94"
$t.Cell(16, 3).Range.Text = "68-67=1"
$t.Cell(16, 4).Range.Text = "64-44=20"
$t.Cell(16, 5).Range.Text = "46-19=27"

$t.Cell(17, 1).Range.Text = "9+15=24"
$t.Cell(17, 2).Range.Text = "32-14=18"
$t.Cell(17, 3).Range.Text = "42-20=22"
$t.Cell(17, 4).Range.Text = "31+6=37"
$t.Cell(17, 5).Range.Text = "61-61=0"

$t.Cell(18, 1).Range.Text = "60-16=44"
$t.Cell(18, 2).Range.Text = "92-1=91"
$t.Cell(18, 3).Range.Text = "85-22=63"
$t.Cell(18, 4).Range.Text = "95-69=26"
$t.Cell(18, 5).Range.Text = "95-48=47"

$t.Cell(19, 1).Range.Text = "17-16=1"
$t.Cell(19, 2).Range.Text = "35+56=91"
$t.Cell(19, 3).Range.Text = "86-72=14"
$t.Cell(19, 4).Range.Text = "52+28=80"
$t.Cell(19, 5).Range.Text = "86-40=46"

$t.Cell(20, 1).Range.Text = "79-59=20"
$t.Cell(20, 2).Range.Text = "9+3=12"
$t.Cell(20, 3).Range.Text = "91-62=29"
$t.Cell(20, 4).Range.Text = "42-15=27"
$t.Cell(20, 5).Range.Text = "82-20=62"

